$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Session" to "Neurology"
$ws.Name = "Neurology"

# Append a new row of scanner log data (row 42). The Student ID column
# holds numeric-looking values that are stored as text (as all the
# other rows in column A already are), so force text formatting before
# writing the value to avoid Excel auto-converting it to a number.
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "201987"
$ws.Range("B42").Value = "Neurology"
$ws.Range("C42").Value = "28/12/2025"
$ws.Range("D42").Value = "10:53:14"
$ws.Range("E42").Value = "Manual"
$ws.Range("F42").Value = "emp17.farah.a.youssef@gmail.com"
